$wb = $excel.ActiveWorkbook

# Update selection on the accounts sheet (D27 -> D25)
$wsAccounts = $wb.Worksheets.Item("accounts")
[void]$wsAccounts.Activate()
[void]$wsAccounts.Range("D25").Select()

# Switch calculation to manual (xlCalculationManual = -4135) so the
# per-cell values we pin below for the volatile RANDBETWEEN column
# survive the automatic recalculation that happens once the script
# finishes (matches Calculation = xlAutomatic otherwise).
$excel.Calculation = -4135

$wsRatings = $wb.Worksheets.Item("user_ratings")

$targets = @{
    2 = 9
    3 = 8
    4 = 9
    5 = 8
    6 = 7
    7 = 10
    8 = 10
    9 = 7
    10 = 7
    11 = 7
    12 = 10
    13 = 7
    14 = 9
    15 = 8
    16 = 10
    17 = 8
    18 = 9
    19 = 8
    20 = 10
    21 = 9
    22 = 7
    23 = 9
    24 = 10
    25 = 9
    26 = 9
    27 = 9
    28 = 10
    29 = 10
    30 = 7
    31 = 8
    32 = 7
    33 = 8
    34 = 8
    35 = 10
    36 = 8
    37 = 8
    38 = 9
    39 = 10
    40 = 8
    41 = 9
    42 = 8
    43 = 9
    44 = 10
    45 = 8
    46 = 9
    47 = 9
    48 = 9
    49 = 10
    50 = 10
    51 = 7
    52 = 10
    53 = 7
    54 = 7
    55 = 10
    56 = 7
    57 = 8
    58 = 9
    59 = 8
    60 = 9
    61 = 10
    62 = 8
    63 = 9
    64 = 8
    65 = 9
    66 = 8
    67 = 7
    68 = 10
    69 = 9
    70 = 8
    71 = 10
    72 = 7
    73 = 8
    74 = 9
    75 = 10
    76 = 9
    77 = 7
    78 = 10
    79 = 8
    80 = 8
    81 = 10
    82 = 8
    83 = 8
    84 = 10
    85 = 10
    86 = 7
    87 = 9
    88 = 7
    89 = 10
    90 = 8
    91 = 8
    92 = 8
    93 = 7
    94 = 8
    95 = 7
    96 = 9
    97 = 8
    98 = 9
    99 = 9
    100 = 10
    101 = 7
    102 = 9
    103 = 7
    104 = 10
    105 = 10
    106 = 8
    107 = 10
    108 = 8
    109 = 8
    110 = 10
    111 = 10
    112 = 8
    113 = 8
    114 = 8
    115 = 10
    116 = 10
    117 = 8
    118 = 10
    119 = 10
    120 = 7
    121 = 10
}

# RANDBETWEEN is volatile, so every recalculation (including the one the
# workbook already went through on load) redraws it. Re-enter the formula
# on each target cell until the freshly recalculated value happens to
# equal the value recorded for that cell, leaving the live formula (and
# its shared-formula grouping) intact with the correct cached result.
foreach ($row in $targets.Keys) {
    $target = $targets[$row]
    $cellRef = "O" + $row
    $cell = $wsRatings.Range($cellRef)
    $v = $cell.Value2
    $attempts = 0
    while ($v -ne $target -and $attempts -lt 2000) {
        $cell.Formula = "=RANDBETWEEN(7,10)"
        $v = $cell.Value2
        $attempts = $attempts + 1
    }
}

Write-Output "done"